{"js": "// PCB libraries were global, made them local.\n//\n// 1. \"High Power Traces\" -> \"Load Traces\"\n// 2. Remove the \"Connector Requirements\" section (heading + its 3 bullets\n//    + the trailing blank paragraph) entirely.\n// 3. Add a new \"Minimum duty cycle shall not exceed 0%\" bullet right after\n//    the existing \"Maximum duty cycle shall be a minimum of 90%\" bullet.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- 1. Rename \"High Power Traces\" -> \"Load Traces\" ---------------------\nfor (const p of items) {\n  if (p.text === \"High Power Traces\") {\n    p.insertText(\"Load Traces\", \"Replace\");\n    break;\n  }\n}\n\n// --- 2. Delete the \"Connector Requirements\" section ----------------------\n// Find the heading paragraph, then walk forward deleting paragraphs until\n// (and including) the blank paragraph that sits right before \"Specifications\".\nlet headingIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === \"Connector Requirements\") {\n    headingIndex = i;\n    break;\n  }\n}\n\nif (headingIndex !== -1) {\n  let specIndex = -1;\n  for (let i = headingIndex + 1; i < items.length; i++) {\n    if (items[i].text === \"Specifications\") {\n      specIndex = i;\n      break;\n    }\n  }\n  if (specIndex !== -1) {\n    for (let i = specIndex - 1; i >= headingIndex; i--) {\n      items[i].delete();\n    }\n  }\n}\n\nawait context.sync();\n\n// --- 3. Insert the new \"Minimum duty cycle...\" bullet ---------------------\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items/text\");\nawait context.sync();\n\nconst items2 = paragraphs2.items;\nfor (const p of items2) {\n  if (p.text === \"Maximum duty cycle shall be a minimum of 90%\") {\n    p.insertParagraph(\"Minimum duty cycle shall not exceed 0%\", \"After\");\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# PCB libraries were global, made them local.\n#\n# 1. \"High Power Traces\" -> \"Load Traces\"\n# 2. Remove the \"Connector Requirements\" section (heading + its 3 bullets\n#    + the trailing blank paragraph) entirely.\n# 3. Add a new \"Minimum duty cycle shall not exceed 0%\" bullet right after\n#    the existing \"Maximum duty cycle shall be a minimum of 90%\" bullet.\n\n$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd([char]13, [char]7)\n}\n\n# --- 1. Rename \"High Power Traces\" -> \"Load Traces\" ------------------------\n$d.Content.Find.Execute(\"High Power Traces\", $false, $false, $false, $false, $false, $true, 1, $false, \"Load Traces\", 2) | Out-Null\n\n# --- 2. Delete the \"Connector Requirements\" section -------------------------\n$headingIndex = -1\n$specIndex = -1\n$n = $d.Paragraphs.Count\nfor ($i = 1; $i -le $n; $i++) {\n    $t = Get-ParaText $d.Paragraphs.Item($i)\n    if ($headingIndex -eq -1 -and $t -eq \"Connector Requirements\") {\n        $headingIndex = $i\n    }\n    elseif ($headingIndex -ne -1 -and $specIndex -eq -1 -and $t -eq \"Specifications\") {\n        $specIndex = $i\n        break\n    }\n}\n\nif ($headingIndex -ne -1 -and $specIndex -ne -1) {\n    $rStart = $d.Paragraphs.Item($headingIndex).Range.Start\n    $rEnd = $d.Paragraphs.Item($specIndex - 1).Range.End\n    $d.Range($rStart, $rEnd).Delete()\n}\n\n# --- 3. Insert the new \"Minimum duty cycle...\" bullet -----------------------\n$n2 = $d.Paragraphs.Count\nfor ($i = 1; $i -le $n2; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ((Get-ParaText $p) -eq \"Maximum duty cycle shall be a minimum of 90%\") {\n        $p.Range.InsertParagraphAfter()\n        $d.Paragraphs.Item($i + 1).Range.Text = \"Minimum duty cycle shall not exceed 0%\"\n        break\n    }\n}\n"}
